$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 113 (item id G113=27775)
$ws.Range("H113").Value = 1665.8334
$ws.Range("I113").Value = 1665.8334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1665.8334
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1588.1666
$ws.Range("N113").ClearContents()
# Row 129 (item id G129=36115)
$ws.Range("H129").Value = 1002.2308
$ws.Range("I129").Value = 798
$ws.Range("J129").Value = 1010.4
$ws.Range("K129").Value = 2394
$ws.Range("L129").Value = 3031.2
$ws.Range("M129").Value = 2606
$ws.Range("N129").Value = -13031.2
# Row 132 (item id G132=44049)
$ws.Range("H132").Value = 192689.8
$ws.Range("I132").Value = 3978.1875
$ws.Range("K132").Value = 11934.5625
$ws.Range("M132").Value = -9404.5625
# Row 138 (item id G138=44169)
$ws.Range("H138").Value = 107781.04
$ws.Range("I138").Value = 2087.7368
$ws.Range("J138").Value = 133861.2
$ws.Range("K138").Value = 6263.2104
$ws.Range("L138").Value = 401583.6
$ws.Range("M138").Value = -1123.2104
$ws.Range("N138").Value = -411863.6

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (item id G32=44147)
$ws.Range("H32").Value = 25280.508
$ws.Range("I32").Value = 18398.723
$ws.Range("J32").Value = 43249.61
$ws.Range("K32").Value = 18398.723
$ws.Range("L32").Value = 43249.61
$ws.Range("M32").Value = -18111.723
$ws.Range("N32").Value = -43823.61
# Row 63 (item id G63=12528)
$ws.Range("H63").Value = 3285.7144
$ws.Range("I63").Value = 3285.7144
$ws.Range("K63").Value = 3285.7144
$ws.Range("M63").Value = -2599.7144
# Row 66 (item id G66=12528)
$ws.Range("H66").Value = 3285.7144
$ws.Range("I66").Value = 3285.7144
$ws.Range("K66").Value = 16428.572
$ws.Range("M66").Value = -12996.572

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 35 (item id G35=2350)
$ws.Range("H35").Value = 32980.5
$ws.Range("J35").Value = 32980.5
$ws.Range("L35").Value = 32980.5
$ws.Range("N35").Value = -33600.5
# Row 82 (item id G82=11877)
$ws.Range("H82").Value = 7730.3125
$ws.Range("I82").Value = 3885.8572
$ws.Range("J82").Value = 34641.5
$ws.Range("K82").Value = 3885.8572
$ws.Range("L82").Value = 34641.5
$ws.Range("M82").Value = -3502.8572
$ws.Range("N82").Value = -35407.5
# Row 85 (item id G85=11877)
$ws.Range("H85").Value = 7730.3125
$ws.Range("I85").Value = 3885.8572
$ws.Range("J85").Value = 34641.5
$ws.Range("K85").Value = 3885.8572
$ws.Range("L85").Value = 34641.5
$ws.Range("M85").Value = -2559.8572
$ws.Range("N85").Value = -37293.5

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4 (item id G4=3742)
$ws.Range("H4").Value = 6785.7144
$ws.Range("J4").Value = 6785.7144
$ws.Range("L4").Value = 6785.7144
$ws.Range("N4").Value = -7009.7144
# Row 58 (item id G58=44021)
$ws.Range("H58").Value = 3521.6667
$ws.Range("I58").Value = 949
$ws.Range("J58").Value = 8667
$ws.Range("K58").Value = 949
$ws.Range("L58").Value = 8667
$ws.Range("M58").Value = -746
$ws.Range("N58").Value = -9073
# Row 99 (item id G99=36198)
$ws.Range("H99").Value = 3536.8235
$ws.Range("I99").Value = 3326.5
$ws.Range("J99").Value = 3723.7778
$ws.Range("K99").Value = 3326.5
$ws.Range("L99").Value = 3723.7778
$ws.Range("M99").Value = -1828.5
$ws.Range("N99").Value = -6719.7778
# Row 126 (item id G126=36198)
$ws.Range("H126").Value = 3536.8235
$ws.Range("I126").Value = 3326.5
$ws.Range("J126").Value = 3723.7778
$ws.Range("K126").Value = 9979.5
$ws.Range("L126").Value = 11171.3334
$ws.Range("M126").Value = -7509.5
$ws.Range("N126").Value = -16111.3334
# Row 136 (item id G136=44021)
$ws.Range("H136").Value = 3521.6667
$ws.Range("I136").Value = 949
$ws.Range("J136").Value = 8667
$ws.Range("K136").Value = 2847
$ws.Range("L136").Value = 26001
$ws.Range("M136").Value = -297
$ws.Range("N136").Value = -31101
# Row 141 (item id G141=43345)
$ws.Range("H141").Value = 74930.12
$ws.Range("J141").Value = 75654.64999999999
$ws.Range("L141").Value = 75654.64999999999
$ws.Range("N141").Value = -86014.64999999999

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4 (item id G4=4650)
$ws.Range("H4").Value = 326200.66
$ws.Range("I4").Value = 918331.5600000001
$ws.Range("J4").Value = 528.65
$ws.Range("K4").Value = 2754994.68
$ws.Range("L4").Value = 1585.95
$ws.Range("M4").Value = -2754882.68
$ws.Range("N4").Value = -1809.95
# Row 25 (item id G25=4709)
$ws.Range("H25").Value = 1499.6666
$ws.Range("J25").Value = 1499.6666
$ws.Range("L25").Value = 4498.9998
$ws.Range("N25").Value = -4836.9998
# Row 30 (item id G30=4709)
$ws.Range("H30").Value = 1499.6666
$ws.Range("J30").Value = 1499.6666
$ws.Range("L30").Value = 4498.9998
$ws.Range("N30").Value = -4702.9998
# Row 104 (item id G104=19807)
$ws.Range("H104").Value = 2792.6667
$ws.Range("J104").Value = 2792.6667
$ws.Range("L104").Value = 8378.000100000001
$ws.Range("N104").Value = -13620.0001
# Row 114 (item id G114=27865)
$ws.Range("H114").Value = 1417.8
$ws.Range("I114").Value = 337.1111
$ws.Range("J114").Value = 2302
$ws.Range("K114").Value = 1011.3333
$ws.Range("L114").Value = 6906
$ws.Range("M114").Value = 2242.6667
$ws.Range("N114").Value = -13414
# Row 117 (item id G117=27870)
$ws.Range("H117").Value = 832.38464
$ws.Range("I117").Value = 306.22223
$ws.Range("K117").Value = 918.66669
$ws.Range("M117").Value = 2523.33331
# Row 131 (item id G131=36060)
$ws.Range("H131").Value = 164794.42
$ws.Range("J131").Value = 176326.84
$ws.Range("L131").Value = 528980.52
$ws.Range("N131").Value = -539060.52

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80 (item id G80=12521)
$ws.Range("H80").Value = 10874.833
$ws.Range("I80").Value = 2937.25
$ws.Range("J80").Value = 26750
$ws.Range("K80").Value = 2937.25
$ws.Range("L80").Value = 26750
$ws.Range("M80").Value = -1939.25
$ws.Range("N80").Value = -28746
# Row 83 (item id G83=12521)
$ws.Range("H83").Value = 10874.833
$ws.Range("I83").Value = 2937.25
$ws.Range("J83").Value = 26750
$ws.Range("K83").Value = 14686.25
$ws.Range("L83").Value = 133750
$ws.Range("M83").Value = -9694.25
$ws.Range("N83").Value = -143734
# Row 102 (item id G102=36169)
$ws.Range("H102").Value = 24254.37
$ws.Range("I102").Value = 16659.445
$ws.Range("J102").Value = 31089.8
$ws.Range("K102").Value = 16659.445
$ws.Range("L102").Value = 31089.8
$ws.Range("M102").Value = -15037.445
$ws.Range("N102").Value = -34333.8
# Row 132 (item id G132=44008)
$ws.Range("H132").Value = 20955.057
$ws.Range("I132").Value = 1762.9736
$ws.Range("K132").Value = 5288.9208
$ws.Range("M132").Value = -2758.9208

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2 (item id G2=2631)
$ws.Range("H2").Value = 25000.666
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 25000.666
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 25000.666
$ws.Range("N2").Value = -25224.666
$ws.Range("M2").ClearContents()
# Row 122 (item id G122=36247)
$ws.Range("H122").Value = 2801.182
$ws.Range("I122").Value = 2973.158
$ws.Range("J122").Value = 2567.7856
$ws.Range("K122").Value = 8919.474
$ws.Range("L122").Value = 7703.3568
$ws.Range("M122").Value = -6469.474
$ws.Range("N122").Value = -12603.3568
# Row 132 (item id G132=44058)
$ws.Range("H132").Value = 187206.39
$ws.Range("I132").Value = 47022.49
$ws.Range("J132").Value = 502620.16
$ws.Range("K132").Value = 141067.47
$ws.Range("L132").Value = 1507860.48
$ws.Range("M132").Value = -138537.47
$ws.Range("N132").Value = -1512920.48
# Row 136 (item id G136=44060)
$ws.Range("H136").Value = 239604.86
$ws.Range("I136").Value = 304283.4
$ws.Range("K136").Value = 912850.2000000001
$ws.Range("M136").Value = -910300.2000000001

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5 (item id G5=3515)
$ws.Range("H5").Value = 39332.332
$ws.Range("J5").Value = 39332.332
$ws.Range("L5").Value = 39332.332
$ws.Range("N5").Value = -39556.332
# Row 136 (item id G136=44031)
$ws.Range("H136").Value = 1538132.2
$ws.Range("I136").Value = 2102624.5
$ws.Range("K136").Value = 6307873.5
$ws.Range("M136").Value = -6305323.5

